$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to remain text-typed like the original inline strings,
# so numeric-looking values (e.g. "313.25", "8.310") are not coerced into numbers
# and trailing zeros / multi-dot "thousand.thousand.decimal" formats are preserved.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.651.56"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.850.11"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").Value = "313.25"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Value = "0.4235"
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("D8").Value = "0.3644"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "44.44"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "0.07293"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "0.8766"
$ws.Range("E11").Value = "  -5.17%  "
$ws.Range("D12").Value = "20.73"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "1.852.58"
$ws.Range("E13").Value = "  -3.57%  "
$ws.Range("D14").Value = "5.343"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "6.527"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").Value = "0.06888"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "79.73"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "0.000008921"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "27.679.11"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").Value = "4.991"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").Value = "10.40"
$ws.Range("E24").Value = "  -4.84%  "
$ws.Range("D25").Value = "2.068.78"
$ws.Range("E25").Value = "  -4.53%  "
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "153.14"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "18.98"
$ws.Range("E28").Value = "  +3.57%  "
$ws.Range("D29").Value = "122.21"
$ws.Range("E29").Value = "  +8.81%  "
$ws.Range("D30").Value = "5.271"
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("D31").Value = "1.887"
$ws.Range("E31").Value = "  +12.74%  "
$ws.Range("D32").Value = "0.08864"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "0.7692"
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("D34").Value = "4.558"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("D35").Value = "2.975"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("D37").Value = "0.9999"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("D39").Value = "0.05357"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").Value = "0.01931"
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").Value = "2.810"
$ws.Range("E41").Value = "  -5.92%  "
$ws.Range("D42").Value = "6.894"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").Value = "0.5107"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").Value = "8.310"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("D46").Value = "0.06545"
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").Value = "0.4768"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "105.68"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "10.31"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").Value = "0.9996"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "1.628"
$ws.Range("E51").Value = "  -1.60%  "

# Restore default styling on column D so no residual number-format style lingers on the cells
$dRange.Style = "Normal"

